$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ---- Block 1: rows 14-21 (audio list update) ----

# Step 1: fill the "#" id column A for the whole first block (A013..A020)
$ws.Cells.Item(14,1).Value = "A013"
$ws.Cells.Item(15,1).Value = "A014"
$ws.Cells.Item(16,1).Value = "A015"
$ws.Cells.Item(17,1).Value = "A016"
$ws.Cells.Item(18,1).Value = "A017"
$ws.Cells.Item(19,1).Value = "A018"
$ws.Cells.Item(20,1).Value = "A019"
$ws.Cells.Item(21,1).Value = "A020"

# Step 2-4: names for rows 14-16
$ws.Cells.Item(14,2).Value = "Attack Shout (human)"
$ws.Cells.Item(15,2).Value = "units march (human)"
$ws.Cells.Item(16,2).Value = "Start Battle"

# Step 5: description for row 16
$ws.Cells.Item(16,7).Value = "world map to combat map, signals beginning of combat, battle drum sounding and slowly fading away in long echo while the view on your screen changes completely from the world map to the combat map in a blurry zoom like effect"

# Step 6-8: names for rows 17-19 (battle music tracks)
$ws.Cells.Item(17,2).Value = "Battle Music Track - Drama 1"
$ws.Cells.Item(18,2).Value = "Battle Music Track - Drama 2"
$ws.Cells.Item(19,2).Value = "Battle Music Track - Drama 3"

# Step 9-11: shared description for rows 17-19
$ws.Cells.Item(17,7).Value = "the music in the combat mode, This music is more intense, faster, and wilder, a strong primeval beat from drums makes your shoulders shrink and your breathing fastens, it increases intensity and drama in multiple layers activated by the remaining hitpoint off all units in the battle, the less HPs remaining activate the next step in intensity and drama"
$ws.Cells.Item(18,7).Value = "the music in the combat mode, This music is more intense, faster, and wilder, a strong primeval beat from drums makes your shoulders shrink and your breathing fastens, it increases intensity and drama in multiple layers activated by the remaining hitpoint off all units in the battle, the less HPs remaining activate the next step in intensity and drama"
$ws.Cells.Item(19,7).Value = "the music in the combat mode, This music is more intense, faster, and wilder, a strong primeval beat from drums makes your shoulders shrink and your breathing fastens, it increases intensity and drama in multiple layers activated by the remaining hitpoint off all units in the battle, the less HPs remaining activate the next step in intensity and drama"

# Step 12-13: names for rows 20-21 (arrow sounds)
$ws.Cells.Item(20,2).Value = "sound of an arrow flying"
$ws.Cells.Item(21,2).Value = "sound of an arrow hitting"

# Step 14: status/category columns for rows 14-21
$ws.Cells.Item(14,3).Value = "SFX"
$ws.Cells.Item(14,4).Value = "open"
$ws.Cells.Item(15,3).Value = "SFX"
$ws.Cells.Item(15,4).Value = "open"
$ws.Cells.Item(16,3).Value = "SFX"
$ws.Cells.Item(16,4).Value = "open"
$ws.Cells.Item(17,3).Value = "soundtrack"
$ws.Cells.Item(17,4).Value = "open"
$ws.Cells.Item(18,3).Value = "soundtrack"
$ws.Cells.Item(18,4).Value = "open"
$ws.Cells.Item(19,3).Value = "soundtrack"
$ws.Cells.Item(19,4).Value = "open"
$ws.Cells.Item(20,3).Value = "SFX"
$ws.Cells.Item(20,4).Value = "open"
$ws.Cells.Item(21,3).Value = "SFX"
$ws.Cells.Item(21,4).Value = "open"

# ---- Block 2: rows 22-29 (particle effects / creature & battle end list update) ----

# Step 15: fill the "#" id column A for the second block (A021..A028)
$ws.Cells.Item(22,1).Value = "A021"
$ws.Cells.Item(23,1).Value = "A022"
$ws.Cells.Item(24,1).Value = "A023"
$ws.Cells.Item(25,1).Value = "A024"
$ws.Cells.Item(26,1).Value = "A025"
$ws.Cells.Item(27,1).Value = "A026"
$ws.Cells.Item(28,1).Value = "A027"
$ws.Cells.Item(29,1).Value = "A028"

# Step 16-21: names for rows 22-27
$ws.Cells.Item(22,2).Value = "Horned Lion roar"
$ws.Cells.Item(23,2).Value = "Horned Lion roar lauder roar"
$ws.Cells.Item(24,2).Value = "unit hit by a claw attack"
$ws.Cells.Item(25,2).Value = "unit shout dying"
$ws.Cells.Item(26,2).Value = "End Battly Vicory"
$ws.Cells.Item(27,2).Value = "End Battly Defeat"

# Step 22: status/category columns for rows 22-27
$ws.Cells.Item(22,3).Value = "SFX"
$ws.Cells.Item(22,4).Value = "open"
$ws.Cells.Item(23,3).Value = "SFX"
$ws.Cells.Item(23,4).Value = "open"
$ws.Cells.Item(24,3).Value = "SFX"
$ws.Cells.Item(24,4).Value = "open"
$ws.Cells.Item(25,3).Value = "SFX"
$ws.Cells.Item(25,4).Value = "open"
$ws.Cells.Item(26,3).Value = "SFX"
$ws.Cells.Item(26,4).Value = "open"
$ws.Cells.Item(27,3).Value = "SFX"
$ws.Cells.Item(27,4).Value = "open"

# Step 23: the combat narrative description for the victory sound (row 26)
$ws.Cells.Item(26,7).Value = "a horn, cheerfull, epic, will give the player feeling of victory, domination, celebration, confidance, rewarding sound"

# ---- final view state: selection resting on B27, scrolled back to show column A ----
$null = $ws.Range("B27").Select()
